$d = $word.ActiveDocument

# The document contains a single tracked change: a deletion of the
# letter "G" in "MA_G08_01_CO" (author "Josué"). Resolve it by accepting
# the deletion, which permanently removes the "G" and leaves the visible
# text as "MA_08_01_CO".
foreach ($rev in $d.Revisions) {
    $rev.Accept()
}

# Word always keeps a single "_GoBack" bookmark marking the location of
# the most recent edit. Locate the end of "MA_" -- i.e. right where the
# just-resolved edit happened -- and move the bookmark there.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("MA_", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

# Remove the old "_GoBack" bookmark whereever it currently sits ...
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ... and recreate it, collapsed, right after "MA_".
$editRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $editRange)
